# Bug fixes for Outlier detection
# Update the adjacency matrix on sheet1: it grows from a 7x7 (A1:G7) matrix
# to a 10x10 matrix (A1:J10), and several existing 0/1 values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final 10x10 adjacency matrix (rows 1-10, columns A-J)
$values = @(
    @(1, 1, 1, 1, 1, 0, 0, 1, 0, 0),
    @(1, 1, 1, 0, 0, 0, 1, 0, 0, 1),
    @(1, 1, 1, 0, 1, 0, 0, 1, 0, 0),
    @(1, 0, 0, 1, 0, 0, 0, 0, 1, 0),
    @(1, 0, 1, 0, 1, 1, 1, 1, 1, 1),
    @(0, 0, 0, 0, 1, 1, 0, 1, 1, 0),
    @(0, 1, 0, 0, 1, 0, 1, 1, 0, 1),
    @(1, 0, 1, 0, 1, 1, 1, 1, 1, 0),
    @(0, 0, 0, 1, 1, 1, 0, 1, 1, 1),
    @(0, 1, 0, 0, 1, 0, 1, 0, 1, 1)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Match the saved selection state from the edit
$ws.Range("I14").Select()
